$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column M (26-jun) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting from L1 (existing last header) onto M1, then set its value.
$ws1.Range("L1").Copy()
$ws1.Range("M1").PasteSpecial(-4122)
$ws1.Range("M1").Value = "26-jun"

$m1Values = @(119.95, 99.63, 93.02, 76.2, 77.45, 79.62, 84.95, 114.4, 106.79, 85, 69.02, 64.44, 62.59, 45.64, 37.39, 36.2, 37.39, 59.13, 80.09, 97.81, 114.53, 110.61, 112.12, 99.84)

for ($i = 0; $i -lt $m1Values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 13).Value = $m1Values[$i]
}

# --- Sheet "Gaz": add row 8 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Cells.Item(8, 1).NumberFormat = "@"
$ws2.Cells.Item(8, 1).Value = "2025-06-25"
$ws2.Cells.Item(8, 1).Style = "Normal"
$ws2.Cells.Item(8, 2).Value = 35.05

# --- Sheet "CO2": add row 8 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Cells.Item(8, 1).NumberFormat = "@"
$ws3.Cells.Item(8, 1).Value = "2025-06-25"
$ws3.Cells.Item(8, 1).Style = "Normal"
$ws3.Cells.Item(8, 2).Value = 70.17
